$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 536 (this shifts existing rows 536-596 down to 537-597,
# preserving their values and formatting, matching the target diff's row-shift pattern).
$ws.Rows.Item(536).Insert()

# Populate the newly inserted (blank) row 536 with a duplicate of the record that is now
# at row 537 (the former row 536), except for the Fecha (D) and Volumen (J) columns, which
# get fresh values for the new weekly data point.
$ws.Cells.Item(536, 1).Value = 6
$ws.Cells.Item(536, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 44946
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = 100112039
$ws.Cells.Item(536, 7).Value = "Ciboulette"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Primera"
$ws.Cells.Item(536, 10).Value = 710
$ws.Cells.Item(536, 11).Value = 900
$ws.Cells.Item(536, 12).Value = 1000
$ws.Cells.Item(536, 13).Value = 954
$ws.Cells.Item(536, 14).Value = "$/docena de atados"
$ws.Cells.Item(536, 15).Value = "Región Metropolitana"
$ws.Cells.Item(536, 16).Value = 318
$ws.Cells.Item(536, 17).Value = 3
$ws.Cells.Item(536, 18).Value = "Hortaliza"
